# role & permission & maps manage. Complete.
# Rows 37-44 (error codes 1035-1042) get re-mapped to their correct
# (previously mismatched) descriptions, and two new rows are appended
# for role/permission management (1043, 1044) ahead of the trailing
# ERROR_TEST (1999) sentinel row, which shifts from row 45 to row 47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 37: 1035 ERROR_NONE_ROLES / authorization / "未发现任何角色"
$ws.Range("A37").Value = 1035
$ws.Range("B37").Value = "ERROR_NONE_ROLES"
$ws.Range("C37").Value = "authorization"
$ws.Range("D37").Value = "未发现任何角色"

# row 38: 1036 ERROR_ROLE_NOT_EXIST / authorization / "指定角色不存在"
$ws.Range("A38").Value = 1036
$ws.Range("B38").Value = "ERROR_ROLE_NOT_EXIST"
$ws.Range("C38").Value = "authorization"
$ws.Range("D38").Value = "指定角色不存在"

# row 39: 1037 ERROR_INVALID_ROLE_NAME / authorization / "用户名无效"
$ws.Range("A39").Value = 1037
$ws.Range("B39").Value = "ERROR_INVALID_ROLE_NAME"
$ws.Range("C39").Value = "authorization"
$ws.Range("D39").Value = "用户名无效"

# row 40: 1038 ERROR_ROLE_EXIST / authorization / "指定角色已存在或重名"
$ws.Range("A40").Value = 1038
$ws.Range("B40").Value = "ERROR_ROLE_EXIST"
$ws.Range("C40").Value = "authorization"
$ws.Range("D40").Value = "指定角色已存在或重名"

# row 41: 1039 ERROR_ADD_ACC_ROLE_FAILED / authorization / "添加账户角色失败"
$ws.Range("A41").Value = 1039
$ws.Range("B41").Value = "ERROR_ADD_ACC_ROLE_FAILED"
$ws.Range("C41").Value = "authorization"
$ws.Range("D41").Value = "添加账户角色失败"

# row 42: 1040 ERROR_NONE_PERMS / authorization / "未发现任何权限"
$ws.Range("A42").Value = 1040
$ws.Range("B42").Value = "ERROR_NONE_PERMS"
$ws.Range("C42").Value = "authorization"
$ws.Range("D42").Value = "未发现任何权限"

# row 43: 1041 ERROR_PERM_NOT_EXIST / authorization / "指定权限不存在"
$ws.Range("A43").Value = 1041
$ws.Range("B43").Value = "ERROR_PERM_NOT_EXIST"
$ws.Range("C43").Value = "authorization"
$ws.Range("D43").Value = "指定权限不存在"

# row 44: 1042 ERROR_PERM_EXIST / authorization / "指定权限已存在或重名"
$ws.Range("A44").Value = 1042
$ws.Range("B44").Value = "ERROR_PERM_EXIST"
$ws.Range("C44").Value = "authorization"
$ws.Range("D44").Value = "指定权限已存在或重名"

# row 45 (new): 1043 ERROR_ADD_ROLE_PERM_FAILED / authorization / "添加角色权限失败"
$ws.Range("A45").Value = 1043
$ws.Range("B45").Value = "ERROR_ADD_ROLE_PERM_FAILED"
$ws.Range("C45").Value = "authorization"
$ws.Range("D45").Value = "添加角色权限失败"

# row 46 (new): 1044 ERROR_DEL_ROLE_PERM_FAILED / authorization / "删除角色权限失败"
$ws.Range("A46").Value = 1044
$ws.Range("B46").Value = "ERROR_DEL_ROLE_PERM_FAILED"
$ws.Range("C46").Value = "authorization"
$ws.Range("D46").Value = "删除角色权限失败"

# row 47: the old trailing ERROR_TEST sentinel row (1999), shifted down from row 45
$ws.Range("A47").Value = 1999
$ws.Range("B47").Value = "ERROR_TEST"
$ws.Range("C47").Value = "general"
$ws.Range("D47").Value = "测试用错误码"

# Update the view state to match the author's final cursor position
# (scroll position + active cell/selection).
$aw = $excel.ActiveWindow
$aw.ScrollRow = 22
$aw.ScrollColumn = 1
$ws.Range("G40").Select() | Out-Null
